# Insert a new weekly price record at row 206 of Sheet1, pushing all
# subsequent rows (old 206..298) down by one (new 207..299). The new
# row 206 re-uses the descriptive/categorical values that the old row
# 206 had (same market/product/category/quality/unit), but carries a
# new date and new volume/price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 206:298 down to 207:299, leaving row 206 free for the
# new record (same effect as right-clicking row 206 -> "Insert").
$ws.Rows.Item(206).Insert()

# Populate the newly inserted row 206 with its data.
$ws.Cells.Item(206, 1).Value  = 7
$ws.Cells.Item(206, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(206, 3).Value  = "Ñuble"
$ws.Cells.Item(206, 4).Value  = 45016
$ws.Cells.Item(206, 5).Value  = 16
$ws.Cells.Item(206, 6).Value  = "Fruta"
$ws.Cells.Item(206, 7).Value  = 100108
$ws.Cells.Item(206, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(206, 9).Value  = 100108005
$ws.Cells.Item(206, 10).Value = "Piña"
$ws.Cells.Item(206, 11).Value = "Caramelo"
$ws.Cells.Item(206, 12).Value = "Segunda"
$ws.Cells.Item(206, 13).Value = 20
$ws.Cells.Item(206, 14).Value = 25000
$ws.Cells.Item(206, 15).Value = 25000
$ws.Cells.Item(206, 16).Value = 25000
$ws.Cells.Item(206, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(206, 18).Value = "Ecuador"
$ws.Cells.Item(206, 19).Value = 1786
$ws.Cells.Item(206, 20).Value = 14

# Make sure the date cell keeps the same date number format that the
# rest of column D uses (style index 2 in the original workbook).
$ws.Cells.Item(206, 4).NumberFormat = $ws.Cells.Item(207, 4).NumberFormat
